$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: fix the username typo and append a missing period to the
# password-mismatch message.
$ws.Range("A7").Value = "The_Big_leauge"
$ws.Range("D7").Value = "password_mismatch:The two password fields didn" + [char]8217 + "t match."

# Move the active selection to D7, matching the saved view state.
$ws.Range("D7").Select()
